$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: updated publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicate "Contact" row; turn it into the "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 duplicated row 10 ("Contact" / "No display for ContactDetail") in the
# original workbook; that duplicate is removed entirely, shifting subsequent
# rows up by one.
$meta.Rows.Item(11).Delete()

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short & Definition now reflect the profile title
# instead of the generic boilerplate text.
$elements.Range("K2").Value = "Insight Confidence Details"
$elements.Range("L2").Value = "Insight Confidence Details"

# The "method" and "score" sub-extensions no longer carry a RIM mapping of
# "N/A" - the mapping cell is cleared.
$elements.Range("AJ5").Value = ""
$elements.Range("AJ6").Value = ""
